$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix default time values in edificios so they don't include quotes
# (previous error). Some additional value tweaks per building as per
# the corrected default data set.

# Anasagasti 1 (row 2)
$ws.Range("B2").Value = "9:00-18:00 "
$ws.Range("C2").Value = "12:00-21:00"
$ws.Range("D2").Value = "9:00-18:00 "
$ws.Range("E2").Value = "12:00-21:00 "
$ws.Range("F2").Value = "9:00-18:00 "
$ws.Range("G2").Value = "9:00-18:00 "

# Anasagasti 2 (row 3)
$ws.Range("B3").Value = "10:00-18:00 "
$ws.Range("C3").Value = "12:00-21:00"
$ws.Range("D3").Value = "9:00-18:00"
$ws.Range("E3").Value = "12:00-21:00 "
$ws.Range("F3").Value = "9:00-18:00 "
$ws.Range("G3").Value = "9:00-12:00 "

# Mitre 1 (row 4)
$ws.Range("B4").Value = "11:00-18:00 "
$ws.Range("C4").Value = "12:00-21:00"
$ws.Range("D4").Value = "9:00-18:00 "
$ws.Range("E4").Value = "12:00-21:00 "
$ws.Range("F4").Value = "9:00-18:00 "
$ws.Range("G4").Value = "9:00-12:00"

$ws.Range("G5").Select()
